$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.591.36'
$ws.Range('E2').Value = '  +2.31%  '
# Row 3
$ws.Range('D3').Value = '1.914.81'
$ws.Range('E3').Value = '  +5.67%  '
# Row 4
$ws.Range('E4').Value = '  +0.00%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.47'
$ws.Range('E5').Value = '  +1.73%  '
# Row 6
$ws.Range('E6').Value = '  -0.07%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5141'
$ws.Range('E7').Value = '  +3.21%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3969'
$ws.Range('E8').Value = '  +1.29%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09818'
$ws.Range('E9').Value = '  -0.91%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.161'
$ws.Range('E10').Value = '  +5.63%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.30'
$ws.Range('E11').Value = '  +3.25%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.551'
$ws.Range('E12').Value = '  +2.43%  '
# Row 13
$ws.Range('E13').Value = '  +4.09%  '
# Row 14
$ws.Range('D14').Value = '1.916.34'
$ws.Range('E14').Value = '  +5.74%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.600'
$ws.Range('E15').Value = '  +4.69%  '
# Row 16
$ws.Range('E16').Value = '  +0.00%  '
# Row 17
$ws.Range('E17').Value = '  +0.14%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.08'
$ws.Range('E18').Value = '  +1.96%  '
# Row 19
$ws.Range('E19').Value = '  +0.17%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.24'
$ws.Range('E20').Value = '  +6.39%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.337'
# Row 23
$ws.Range('D23').Value = '28.643.10'
$ws.Range('E23').Value = '  +2.28%  '
# Row 24
$ws.Range('E24').Value = '  +3.68%  '
# Row 25
$ws.Range('E25').Value = '  +1.29%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.708'
$ws.Range('E26').Value = '  +13.27%  '
# Row 27
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.395'
$ws.Range('E27').Value = '  -0.81%  '
# Row 28
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.133.16'
$ws.Range('E28').Value = '  +5.45%  '
# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '21.29'
$ws.Range('E29').Value = '  +3.58%  '
# Row 30
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '160.03'
$ws.Range('E30').Value = '  +1.02%  '
# Row 31
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '129.05'
$ws.Range('E31').Value = '  +2.04%  '
# Row 32
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.107'
$ws.Range('E32').Value = '  +7.30%  '
# Row 33
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1083'
$ws.Range('E33').Value = '  +1.81%  '
# Row 34
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.768'
$ws.Range('E34').Value = '  +3.85%  '
# Row 35
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.640'
$ws.Range('E35').Value = '  +1.21%  '
# Row 36
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.874'
$ws.Range('E36').Value = '  +11.03%  '
# Row 37
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06822'
$ws.Range('E37').Value = '  +1.61%  '
# Row 38
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02442'
$ws.Range('E38').Value = '  +5.04%  '
# Row 39
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.271'
$ws.Range('E39').Value = '  +8.13%  '
# Row 40
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2238'
$ws.Range('E40').Value = '  +4.57%  '
# Row 41
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.93'
$ws.Range('E41').Value = '  +6.11%  '
# Row 42
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.131'
$ws.Range('E42').Value = '  +4.15%  '
# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6454'
$ws.Range('E43').Value = '  +4.52%  '
# Row 44
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.195'
$ws.Range('E44').Value = '  +2.15%  '
# Row 45
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9998'
$ws.Range('E45').Value = '  -0.14%  '
# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.65'
$ws.Range('E46').Value = '  +3.51%  '
# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6101'
$ws.Range('E47').Value = '  +3.59%  '
# Row 48
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.793'
$ws.Range('E48').Value = '  +2.79%  '
# Row 49
$ws.Range('B49').Value = 'WEMIXTOKEN'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.279'
$ws.Range('E49').Value = '  -0.19%  '
# Row 50
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.041'
$ws.Range('E50').Value = '  +5.69%  '
# Row 51
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.53'
$ws.Range('E51').Value = '  +1.61%  '
